# Applies the "riska.xlsx" update:
#  - rename sheet (3) -> (4)
#  - refresh several Talk_time / Repayment_collections metrics
#  - update Fadilah Damayanti's repayment amount & recovery rate
#  - fill in Yandi Nugraha's previously-blank repayment figures

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- rename worksheet tab ---------------------------------------------
$ws.Name = "repayment_20250912_20250912 (4)"

# --- helper: write a literal text value (keeps shared-string type and
#     avoids leaving a lingering number-format style on the cell) ------
function Set-TextCell($rng, [string]$text) {
    $rng.NumberFormat = "@"
    $rng.Value = $text
    $rng.Style = "Normal"
}

# --- plain numeric metric refreshes ------------------------------------
$ws.Range("H2").Value = 248
$ws.Range("H3").Value = 1.574
$ws.Range("D4").Value = 6
$ws.Range("H4").Value = 1.642
$ws.Range("H5").Value = 1.373
$ws.Range("H6").Value = 1.123
$ws.Range("H7").Value = 1.522
$ws.Range("H8").Value = 1.061
$ws.Range("H9").Value = 637
$ws.Range("H10").Value = 1.027
$ws.Range("H11").Value = 1.152
$ws.Range("H12").Value = 1.209
$ws.Range("H13").Value = 619
$ws.Range("H14").Value = 994
$ws.Range("H16").Value = 883
$ws.Range("H17").Value = 1.393
$ws.Range("H18").Value = 737

# --- Fadilah Damayanti (row 4): repayment amount & recovery rate ------
Set-TextCell $ws.Range("E4") "2,827,991.00"
Set-TextCell $ws.Range("G4") "1.96"

# --- Yandi Nugraha (row 15): was blank ("0.00"), now populated --------
$ws.Range("D15").Value = 2
Set-TextCell $ws.Range("E15") "1,937,673.00"
Set-TextCell $ws.Range("G15") "1.27"
$ws.Range("H15").Value = 4.81
